# Outstandings.xlsx — "Add files via upload" edit
#
# Net effect (per the target diff): on the "Purchase 22-23" sheet the
# entry that used to sit in row 11 (Invoice "PAN9562/23-24" / vendor
# "Microciti", amount 1145) was removed, which shifts every later
# entry up by one slot (two worksheet rows, since each entry occupies
# one data row followed by a blank spacer row). The last entry
# (old row 21, "3054" / "Bale & Sons", amount 25620) falls off the
# bottom and disappears entirely. The "Sr. No" column is renumbered
# to stay sequential (1..7) after the deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Purchase 22-23")

# Remove the old row-11 entry together with its blank spacer row
# (row 12). Everything below (old rows 13-21) shifts up by two rows,
# carrying its values/formulas/number formats/fonts with it — which is
# exactly why the new row 11 picks up the taller (14.4pt) row height
# that used to belong to old row 13, and new row 13 reverts to the
# default height that used to belong to old row 15, and so on.
$ws.Rows("11:12").Delete()

# The very last entry (old row 21, now shifted to row 19-20) is also
# removed — it no longer appears anywhere in the sheet.
$ws.Rows("19:20").Delete()

# Renumber the "Sr. No" column (col A) for the rows that shifted up.
$ws.Range("A11").Value = 4
$ws.Range("A13").Value = 5
$ws.Range("A15").Value = 6
$ws.Range("A17").Value = 7

# Restore the view state recorded for the sheet after the edit: the
# window was scrolled down a couple of rows and the last selection
# made was on B23 (below the now-shorter used range).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B23").Select()
